$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("2023-12-13 12:07:58", 0.0004),
    @("2023-12-13 12:08:29", 0.0018),
    @("2023-12-13 12:09:21", 0.0036),
    @("2023-12-13 12:09:26", 0.0004),
    @("2023-12-13 12:09:37", 0.0002)
)

$startRow = 273
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}
